# Update the Markov transition-probability matrix on Sheet1 with the
# refreshed figures from the larger simulation run (added more games,
# sped up simulate game logic, and drafted optimization logic).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1554054054054054
$ws.Range("C2").Value = 0.6148648648648649
$ws.Range("J2").Value = 0.01351351351351351
$ws.Range("P2").Value = 0.1182432432432432
$ws.Range("S2").Value = 0.09797297297297297
$ws.Range("B3").Value = 0.02105263157894737
$ws.Range("C3").Value = 0.05263157894736842
$ws.Range("J3").Value = 0.03157894736842105
$ws.Range("P3").Value = 0.7105263157894737
$ws.Range("S3").Value = 0.1842105263157895
$ws.Range("J4").Value = 0.02380952380952381
$ws.Range("P4").Value = 0.7142857142857143
$ws.Range("S4").Value = 0.2619047619047619
$ws.Range("B6").Value = 0.07327586206896551
$ws.Range("D6").Value = 0.004310344827586207
$ws.Range("F6").Value = 0.0603448275862069
$ws.Range("J6").Value = 0.2758620689655172
$ws.Range("O6").Value = 0.03017241379310345
$ws.Range("Q6").Value = 0.1681034482758621
$ws.Range("R6").Value = 0.0603448275862069
$ws.Range("S6").Value = 0.3275862068965517
$ws.Range("B7").Value = 0.09606986899563319
$ws.Range("D7").Value = 0.01310043668122271
$ws.Range("F7").Value = 0.08733624454148471
$ws.Range("J7").Value = 0.08733624454148471
$ws.Range("O7").Value = 0.01746724890829694
$ws.Range("Q7").Value = 0.1921397379912664
$ws.Range("R7").Value = 0.0611353711790393
$ws.Range("S7").Value = 0.445414847161572
$ws.Range("B8").Value = 0.1048158640226629
$ws.Range("D8").Value = 0.0169971671388102
$ws.Range("E8").Value = 0.0028328611898017
$ws.Range("F8").Value = 0.04815864022662889
$ws.Range("J8").Value = 0.141643059490085
$ws.Range("O8").Value = 0.0141643059490085
$ws.Range("Q8").Value = 0.1784702549575071
$ws.Range("R8").Value = 0.1048158640226629
$ws.Range("S8").Value = 0.3881019830028329
$ws.Range("B9").Value = 0.09545454545454546
$ws.Range("D9").Value = 0.004545454545454545
$ws.Range("F9").Value = 0.1045454545454545
$ws.Range("J9").Value = 0.1454545454545454
$ws.Range("O9").Value = 0.00909090909090909
$ws.Range("Q9").Value = 0.1727272727272727
$ws.Range("R9").Value = 0.09545454545454546
$ws.Range("S9").Value = 0.3727272727272727
$ws.Range("B10").Value = 0.1106161841128434
$ws.Range("D10").Value = 0.0244988864142539
$ws.Range("E10").Value = 0.001484780994803266
$ws.Range("F10").Value = 0.066815144766147
$ws.Range("J10").Value = 0.1343726800296956
$ws.Range("O10").Value = 0.01410541945063103
$ws.Range("Q10").Value = 0.2249443207126949
$ws.Range("R10").Value = 0.08537490720118783
$ws.Range("S10").Value = 0.3377876763177431
$ws.Range("G11").Value = 0.1476323119777159
$ws.Range("J11").Value = 0.07799442896935933
$ws.Range("K11").Value = 0.2005571030640668
$ws.Range("L11").Value = 0.5543175487465181
$ws.Range("S11").Value = 0.01949860724233983
$ws.Range("G12").Value = 0.7537688442211056
$ws.Range("J12").Value = 0.2060301507537688
$ws.Range("K12").Value = 0.01005025125628141
$ws.Range("L12").Value = 0.01507537688442211
$ws.Range("S12").Value = 0.01507537688442211
$ws.Range("G13").Value = 0.5833333333333334
$ws.Range("J13").Value = 0.3541666666666667
$ws.Range("S13").Value = 0.0625
$ws.Range("F15").Value = 0.03902439024390244
$ws.Range("H15").Value = 0.1170731707317073
$ws.Range("I15").Value = 0.0975609756097561
$ws.Range("J15").Value = 0.3609756097560975
$ws.Range("K15").Value = 0.05365853658536585
$ws.Range("M15").Value = 0.02439024390243903
$ws.Range("O15").Value = 0.05365853658536585
$ws.Range("S15").Value = 0.2536585365853659
$ws.Range("F16").Value = 0.03092783505154639
$ws.Range("H16").Value = 0.1288659793814433
$ws.Range("I16").Value = 0.1030927835051546
$ws.Range("J16").Value = 0.4278350515463917
$ws.Range("K16").Value = 0.1134020618556701
$ws.Range("M16").Value = 0.03608247422680412
$ws.Range("O16").Value = 0.02061855670103093
$ws.Range("S16").Value = 0.1391752577319588
$ws.Range("F17").Value = 0.01232032854209446
$ws.Range("H17").Value = 0.1273100616016427
$ws.Range("I17").Value = 0.09856262833675565
$ws.Range("J17").Value = 0.4394250513347023
$ws.Range("K17").Value = 0.108829568788501
$ws.Range("M17").Value = 0.01848049281314168
$ws.Range("N17").Value = 0.002053388090349076
$ws.Range("O17").Value = 0.08008213552361396
$ws.Range("S17").Value = 0.1129363449691992
$ws.Range("F18").Value = 0.01507537688442211
$ws.Range("H18").Value = 0.1005025125628141
$ws.Range("I18").Value = 0.1055276381909548
$ws.Range("J18").Value = 0.4371859296482412
$ws.Range("K18").Value = 0.1608040201005025
$ws.Range("M18").Value = 0.005025125628140704
$ws.Range("N18").Value = 0.005025125628140704
$ws.Range("O18").Value = 0.07035175879396985
$ws.Range("S18").Value = 0.1005025125628141
$ws.Range("F19").Value = 0.01380991064175467
$ws.Range("H19").Value = 0.1852152721364744
$ws.Range("I19").Value = 0.0901705930138099
$ws.Range("J19").Value = 0.3744922826969943
$ws.Range("K19").Value = 0.1299756295694557
$ws.Range("M19").Value = 0.02193338748984565
$ws.Range("O19").Value = 0.06173842404549147
$ws.Range("S19").Value = 0.1226645004061738
